# Update New Orleans xlsx file:
#  1. hotel_info: insert a new "State" column (with value "Louisiana" for
#     the single data row) right after "Hotel_Name" and before "City".
#  2. Reorder the worksheets so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

# --- 1. Add the State column to hotel_info -------------------------------
$hotelInfo = $wb.Worksheets.Item("hotel_info")

# Hotel_Name is column B, City is column C -> insert a new column at C.
$hotelInfo.Columns.Item(3).Insert()
$hotelInfo.Range("C1").Value = "State"
$hotelInfo.Range("C2").Value = "Louisiana"

# --- 2. Reorder sheets: review_info first, hotel_info second -------------
$reviewInfo = $wb.Worksheets.Item("review_info")
$reviewInfo.Move($hotelInfo)
